$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): reorder to Username, Password, First Name, Last Name, Age ---
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "First Name"
$ws.Range("D1").Value = "Last Name"
$ws.Range("E1").Value = "Age"

# --- Row 2: replace with the (corrected) single remaining data row ---
# Force the Age column to be stored as text (matches original convention where ages were text)
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "Mezix"
$ws.Range("B2").Value = "password"
$ws.Range("C2").Value = "Felix"
$ws.Range("D2").Value = "Swimmer"
$ws.Range("E2").Value = "21"

# Row 2 no longer carries the centered/custom formatting it used to have
$ws.Range("A2:E2").Style = "Normal"

# --- Row 3: the old duplicate row is removed entirely ---
$ws.Rows.Item(3).Delete()

# --- Column widths updated for the new column order/content ---
$ws.Columns.Item(3).ColumnWidth = 30
$ws.Columns.Item(4).ColumnWidth = 30
$ws.Columns.Item(5).ColumnWidth = 8

# --- Selection moves to C2 ---
$ws.Range("C2").Select() | Out-Null

Write-Host "edit complete"
